$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 15.84401633333333
$ws.Cells.Item(2, 8).Value = 47.532049
$ws.Cells.Item(2, 9).Value = 0.7212145038223592
$ws.Cells.Item(2, 10).Value = 0.7256686596529557
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.808665
$ws.Cells.Item(2, 14).Value = 5.425995
$ws.Cells.Item(2, 15).Value = 0.01261203239065773
$ws.Cells.Item(2, 16).Value = 0.01316260230515319
$ws.Cells.Item(2, 17).Value = 28.65651780152834
$ws.Cells.Item(2, 18).Value = 257.908660213755
$ws.Cells.Item(2, 19).Value = 0.00909598068281974
$ws.Cells.Item(2, 20).Value = 0.009551687972325421

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 15.84401633333333
$ws.Cells.Item(3, 8).Value = 47.532049
$ws.Cells.Item(3, 9).Value = 0.7212145038223592
$ws.Cells.Item(3, 10).Value = 0.7256686596529557
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 21.254561
$ws.Cells.Item(3, 14).Value = 63.763683
$ws.Cells.Item(3, 15).Value = 0.1482105374854993
$ws.Cells.Item(3, 16).Value = 0.1546805702623864
$ws.Cells.Item(3, 17).Value = 336.7576116418297
$ws.Cells.Item(3, 18).Value = 3030.818504776467
$ws.Cells.Item(3, 19).Value = 0.1068915892538495
$ws.Cells.Item(3, 20).Value = 0.1122468420966608

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 15.84401633333333
$ws.Cells.Item(4, 8).Value = 47.532049
$ws.Cells.Item(4, 9).Value = 0.7212145038223592
$ws.Cells.Item(4, 10).Value = 0.7256686596529557
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 46.05851666666666
$ws.Cells.Item(4, 14).Value = 138.17555
$ws.Cells.Item(4, 15).Value = 0.3211714187346186
$ws.Cells.Item(4, 16).Value = 0.3351919441403484
$ws.Cells.Item(4, 17).Value = 729.7518903557722
$ws.Cells.Item(4, 18).Value = 6567.767013201949
$ws.Cells.Item(4, 19).Value = 0.2316334854046111
$ws.Cells.Item(4, 20).Value = 0.243238288830795

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 15.84401633333333
$ws.Cells.Item(5, 8).Value = 47.532049
$ws.Cells.Item(5, 9).Value = 0.7212145038223592
$ws.Cells.Item(5, 10).Value = 0.7256686596529557
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 56.290605
$ws.Cells.Item(5, 14).Value = 168.871815
$ws.Cells.Item(5, 15).Value = 0.3925209663203081
$ws.Cells.Item(5, 16).Value = 0.4096562089339196
$ws.Cells.Item(5, 17).Value = 891.869265033215
$ws.Cells.Item(5, 18).Value = 8026.823385298935
$ws.Cells.Item(5, 19).Value = 0.283091813964574
$ws.Cells.Item(5, 20).Value = 0.2972746720555886

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 15.84401633333333
$ws.Cells.Item(6, 8).Value = 47.532049
$ws.Cells.Item(6, 9).Value = 0.7212145038223592
$ws.Cells.Item(6, 10).Value = 0.7256686596529557
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 17.995546
$ws.Cells.Item(6, 14).Value = 35.991092
$ws.Cells.Item(6, 15).Value = 0.1254850450689161
$ws.Cells.Item(6, 16).Value = 0.08730867435819248
$ws.Cells.Item(6, 17).Value = 285.1217247512513
$ws.Cells.Item(6, 18).Value = 1710.730348507508
$ws.Cells.Item(6, 19).Value = 0.09050163451650474
$ws.Cells.Item(6, 20).Value = 0.06335716869758591

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.951312666666666
$ws.Cells.Item(7, 8).Value = 8.853938
$ws.Cells.Item(7, 9).Value = 0.1343427989301267
$ws.Cells.Item(7, 10).Value = 0.1351724879588164
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.808665
$ws.Cells.Item(7, 14).Value = 5.425995
$ws.Cells.Item(7, 15).Value = 0.01261203239065773
$ws.Cells.Item(7, 16).Value = 0.01316260230515319
$ws.Cells.Item(7, 17).Value = 5.337935924256667
$ws.Cells.Item(7, 18).Value = 48.04142331831
$ws.Cells.Item(7, 19).Value = 0.001694335731558377
$ws.Cells.Item(7, 20).Value = 0.001779221701600009

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.951312666666666
$ws.Cells.Item(8, 8).Value = 8.853938
$ws.Cells.Item(8, 9).Value = 0.1343427989301267
$ws.Cells.Item(8, 10).Value = 0.1351724879588164
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 21.254561
$ws.Cells.Item(8, 14).Value = 63.763683
$ws.Cells.Item(8, 15).Value = 0.1482105374854993
$ws.Cells.Item(8, 16).Value = 0.1546805702623864
$ws.Cells.Item(8, 17).Value = 62.72885510373933
$ws.Cells.Item(8, 18).Value = 564.559695933654
$ws.Cells.Item(8, 19).Value = 0.01991101843674044
$ws.Cells.Item(8, 20).Value = 0.02090855752125528

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2.951312666666666
$ws.Cells.Item(9, 8).Value = 8.853938
$ws.Cells.Item(9, 9).Value = 0.1343427989301267
$ws.Cells.Item(9, 10).Value = 0.1351724879588164
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 46.05851666666666
$ws.Cells.Item(9, 14).Value = 138.17555
$ws.Cells.Item(9, 15).Value = 0.3211714187346186
$ws.Cells.Item(9, 16).Value = 0.3351919441403484
$ws.Cells.Item(9, 17).Value = 135.9330836462111
$ws.Cells.Item(9, 18).Value = 1223.3977528159
$ws.Cells.Item(9, 19).Value = 0.04314706732916841
$ws.Cells.Item(9, 20).Value = 0.0453087290332035

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 2.951312666666666
$ws.Cells.Item(10, 8).Value = 8.853938
$ws.Cells.Item(10, 9).Value = 0.1343427989301267
$ws.Cells.Item(10, 10).Value = 0.1351724879588164
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 56.290605
$ws.Cells.Item(10, 14).Value = 168.871815
$ws.Cells.Item(10, 15).Value = 0.3925209663203081
$ws.Cells.Item(10, 16).Value = 0.4096562089339196
$ws.Cells.Item(10, 17).Value = 166.13117555083
$ws.Cells.Item(10, 18).Value = 1495.18057995747
$ws.Cells.Item(10, 19).Value = 0.0527323652542282
$ws.Cells.Item(10, 20).Value = 0.05537424896937462

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 2.951312666666666
$ws.Cells.Item(11, 8).Value = 8.853938
$ws.Cells.Item(11, 9).Value = 0.1343427989301267
$ws.Cells.Item(11, 10).Value = 0.1351724879588164
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 17.995546
$ws.Cells.Item(11, 14).Value = 35.991092
$ws.Cells.Item(11, 15).Value = 0.1254850450689161
$ws.Cells.Item(11, 16).Value = 0.08730867435819248
$ws.Cells.Item(11, 17).Value = 53.11048285338266
$ws.Cells.Item(11, 18).Value = 318.662897120296
$ws.Cells.Item(11, 19).Value = 0.01685801217843129
$ws.Cells.Item(11, 20).Value = 0.01180173073338299

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 2.768664
$ws.Cells.Item(12, 8).Value = 8.305992
$ws.Cells.Item(12, 9).Value = 0.1260286906426543
$ws.Cells.Item(12, 10).Value = 0.1268070324872419
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 1.808665
$ws.Cells.Item(12, 14).Value = 5.425995
$ws.Cells.Item(12, 15).Value = 0.01261203239065773
$ws.Cells.Item(12, 16).Value = 0.01316260230515319
$ws.Cells.Item(12, 17).Value = 5.00758567356
$ws.Cells.Item(12, 18).Value = 45.06827106204
$ws.Cells.Item(12, 19).Value = 0.001589477928537339
$ws.Cells.Item(12, 20).Value = 0.001669110538126206

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 2.768664
$ws.Cells.Item(13, 8).Value = 8.305992
$ws.Cells.Item(13, 9).Value = 0.1260286906426543
$ws.Cells.Item(13, 10).Value = 0.1268070324872419
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 21.254561
$ws.Cells.Item(13, 14).Value = 63.763683
$ws.Cells.Item(13, 15).Value = 0.1482105374854993
$ws.Cells.Item(13, 16).Value = 0.1546805702623864
$ws.Cells.Item(13, 17).Value = 58.84673787650399
$ws.Cells.Item(13, 18).Value = 529.620640888536
$ws.Cells.Item(13, 19).Value = 0.0186787799787415
$ws.Cells.Item(13, 20).Value = 0.01961458409840753

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 2.768664
$ws.Cells.Item(14, 8).Value = 8.305992
$ws.Cells.Item(14, 9).Value = 0.1260286906426543
$ws.Cells.Item(14, 10).Value = 0.1268070324872419
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 46.05851666666666
$ws.Cells.Item(14, 14).Value = 138.17555
$ws.Cells.Item(14, 15).Value = 0.3211714187346186
$ws.Cells.Item(14, 16).Value = 0.3351919441403484
$ws.Cells.Item(14, 17).Value = 127.5205569884
$ws.Cells.Item(14, 18).Value = 1147.6850128956
$ws.Cells.Item(14, 19).Value = 0.04047681337496763
$ws.Cells.Item(14, 20).Value = 0.04250469575006693

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 2.768664
$ws.Cells.Item(15, 8).Value = 8.305992
$ws.Cells.Item(15, 9).Value = 0.1260286906426543
$ws.Cells.Item(15, 10).Value = 0.1268070324872419
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 56.290605
$ws.Cells.Item(15, 14).Value = 168.871815
$ws.Cells.Item(15, 15).Value = 0.3925209663203081
$ws.Cells.Item(15, 16).Value = 0.4096562089339196
$ws.Cells.Item(15, 17).Value = 155.84977160172
$ws.Cells.Item(15, 18).Value = 1402.64794441548
$ws.Cells.Item(15, 19).Value = 0.04946890343513783
$ws.Cells.Item(15, 20).Value = 0.05194728819488389

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 2.768664
$ws.Cells.Item(16, 8).Value = 8.305992
$ws.Cells.Item(16, 9).Value = 0.1260286906426543
$ws.Cells.Item(16, 10).Value = 0.1268070324872419
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 17.995546
$ws.Cells.Item(16, 14).Value = 35.991092
$ws.Cells.Item(16, 15).Value = 0.1254850450689161
$ws.Cells.Item(16, 16).Value = 0.08730867435819248
$ws.Cells.Item(16, 17).Value = 49.823620370544
$ws.Cells.Item(16, 18).Value = 298.941722223264
$ws.Cells.Item(16, 19).Value = 0.01581471592526996
$ws.Cells.Item(16, 20).Value = 0.01107135390575733

$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.4045285
$ws.Cells.Item(17, 8).Value = 0.809057
$ws.Cells.Item(17, 9).Value = 0.01841400660485959
$ws.Cells.Item(17, 10).Value = 0.01235181990098599
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 1.808665
$ws.Cells.Item(17, 14).Value = 5.425995
$ws.Cells.Item(17, 15).Value = 0.01261203239065773
$ws.Cells.Item(17, 16).Value = 0.01316260230515319
$ws.Cells.Item(17, 17).Value = 0.7316565394525001
$ws.Cells.Item(17, 18).Value = 4.389939236715001
$ws.Cells.Item(17, 19).Value = 0.0002322380477422746
$ws.Cells.Item(17, 20).Value = 0.0001625820931015553

$ws.Cells.Item(18, 5).Value = 2
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 0.4045285
$ws.Cells.Item(18, 8).Value = 0.809057
$ws.Cells.Item(18, 9).Value = 0.01841400660485959
$ws.Cells.Item(18, 10).Value = 0.01235181990098599
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 21.254561
$ws.Cells.Item(18, 14).Value = 63.763683
$ws.Cells.Item(18, 15).Value = 0.1482105374854993
$ws.Cells.Item(18, 16).Value = 0.1546805702623864
$ws.Cells.Item(18, 17).Value = 8.5980756794885
$ws.Cells.Item(18, 18).Value = 51.588454076931
$ws.Cells.Item(18, 19).Value = 0.002729149816167774
$ws.Cells.Item(18, 20).Value = 0.001910586546062807

$ws.Cells.Item(19, 5).Value = 2
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 0.4045285
$ws.Cells.Item(19, 8).Value = 0.809057
$ws.Cells.Item(19, 9).Value = 0.01841400660485959
$ws.Cells.Item(19, 10).Value = 0.01235181990098599
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 46.05851666666666
$ws.Cells.Item(19, 14).Value = 138.17555
$ws.Cells.Item(19, 15).Value = 0.3211714187346186
$ws.Cells.Item(19, 16).Value = 0.3351919441403484
$ws.Cells.Item(19, 17).Value = 18.63198265939166
$ws.Cells.Item(19, 18).Value = 111.79189595635
$ws.Cells.Item(19, 19).Value = 0.005914052625871392
$ws.Cells.Item(19, 20).Value = 0.004140230526282942

$ws.Cells.Item(20, 5).Value = 2
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 0.4045285
$ws.Cells.Item(20, 8).Value = 0.809057
$ws.Cells.Item(20, 9).Value = 0.01841400660485959
$ws.Cells.Item(20, 10).Value = 0.01235181990098599
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 12).Value = 1
$ws.Cells.Item(20, 13).Value = 56.290605
$ws.Cells.Item(20, 14).Value = 168.871815
$ws.Cells.Item(20, 15).Value = 0.3925209663203081
$ws.Cells.Item(20, 16).Value = 0.4096562089339196
$ws.Cells.Item(20, 17).Value = 22.7711540047425
$ws.Cells.Item(20, 18).Value = 136.626924028455
$ws.Cells.Item(20, 19).Value = 0.007227883666368022
$ws.Cells.Item(20, 20).Value = 0.005059999714072465

$ws.Cells.Item(21, 5).Value = 2
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 0.4045285
$ws.Cells.Item(21, 8).Value = 0.809057
$ws.Cells.Item(21, 9).Value = 0.01841400660485959
$ws.Cells.Item(21, 10).Value = 0.01235181990098599
$ws.Cells.Item(21, 11).Value = 2
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = 17.995546
$ws.Cells.Item(21, 14).Value = 35.991092
$ws.Cells.Item(21, 15).Value = 0.1254850450689161
$ws.Cells.Item(21, 16).Value = 0.08730867435819248
$ws.Cells.Item(21, 17).Value = 7.279711230061
$ws.Cells.Item(21, 18).Value = 29.118844920244
$ws.Cells.Item(21, 19).Value = 0.002310682448710125
$ws.Cells.Item(21, 20).Value = 0.001078421021466227

